# Insert a new data row at sheet row 105 (pushes existing rows 105-155 down
# to 106-156) and populate it with the new record, per the commit's weekly
# price-update diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 105..155 down by one, creating a fresh (blank) row 105 that
# inherits formatting (incl. the date number-format on column D) from the
# row above it - matching native Excel "insert row" behaviour.
$ws.Rows.Item(105).Insert()

# Populate the newly inserted row 105 with the new observation.
$ws.Cells.Item(105, 1).Value  = 9
$ws.Cells.Item(105, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(105, 3).Value  = "Metropolitana"
$ws.Cells.Item(105, 4).Value  = 45007
$ws.Cells.Item(105, 5).Value  = 13
$ws.Cells.Item(105, 6).Value  = 100112022
$ws.Cells.Item(105, 7).Value  = "Arveja Verde"
$ws.Cells.Item(105, 8).Value  = "Perfection"
$ws.Cells.Item(105, 9).Value  = "Primera"
$ws.Cells.Item(105, 10).Value = 43
$ws.Cells.Item(105, 11).Value = 27000
$ws.Cells.Item(105, 12).Value = 29000
$ws.Cells.Item(105, 13).Value = 28023
$ws.Cells.Item(105, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(105, 15).Value = "Provincia de Huasco"
$ws.Cells.Item(105, 16).Value = 1121
$ws.Cells.Item(105, 17).Value = 25
$ws.Cells.Item(105, 18).Value = "Hortaliza"
